$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.021.55'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = '3.158.45'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = "'579.73"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +1.68%  '
$ws.Range("D6").Value = "'149.85"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.154.53'
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("E9").Value = '  -0.23%  '
$ws.Range("E10").Value = '  -1.92%  '
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("E12").Value = '  -0.50%  '
$ws.Range("D13").Value = "'0.0000264"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  +2.01%  '
$ws.Range("D14").Value = "'37.22"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  -2.32%  '
$ws.Range("D15").Value = '3.672.40'
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("D16").Value = '64.905.92'
$ws.Range("E16").Value = '  +0.17%  '
$ws.Range("D17").Value = '3.154.33'
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("E18").Value = '  -0.81%  '
$ws.Range("E19").Value = '  +0.54%  '
$ws.Range("D20").Value = "'505.58"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -1.80%  '
$ws.Range("D21").Value = "'14.93"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("D22").Value = "'0.716"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -2.59%  '
$ws.Range("D23").Value = "'15.19"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -0.72%  '
$ws.Range("D24").Value = "'7.74"
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -1.34%  '
$ws.Range("D25").Value = "'84.57"
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  -0.58%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").Value = "'9.04"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +2.88%  '
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("D29").Value = "'2.19"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  +4.69%  '
$ws.Range("D31").Value = "'27.63"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -0.89%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("E33").Value = '  +0.54%  '
$ws.Range("D34").Value = "'6.43"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +4.50%  '
$ws.Range("D35").Value = "'6.49"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -1.54%  '
$ws.Range("D36").Value = "'54.81"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -1.63%  '
$ws.Range("E37").Value = '  +3.04%  '
$ws.Range("D38").Value = "'479.19"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -1.19%  '
$ws.Range("E39").Value = '  -1.96%  '
$ws.Range("D40").Value = "'2.93"
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  -1.33%  '
$ws.Range("E41").Value = '  +0.93%  '
$ws.Range("D42").Value = '2.992.69'
$ws.Range("E42").Value = '  -3.98%  '
$ws.Range("E43").Value = '  -2.01%  '
$ws.Range("D44").Value = "'0.283"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  -3.53%  '
$ws.Range("E45").Value = '  -1.40%  '
$ws.Range("D46").Value = "'28.34"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -3.93%  '
$ws.Range("D47").Value = '0.0₃0593'
$ws.Range("E47").Value = '  +3.23%  '
$ws.Range("D48").Value = "'1.00"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("E49").Value = '  -1.10%  '
$ws.Range("E50").Value = '  -2.38%  '
$ws.Range("D51").Value = "'2.48"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +14.04%  '
